# Reorders the "Recorded By" (column G) comma-separated list of names/emails
# so that any token equal to "System" (case-insensitive) is moved to the
# front of the list, preserving the relative order of the remaining tokens.
# Cells where "System" is already first, or where no "System" token exists,
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $tokens = $val -split ", "

    $systemTokens = @()
    $otherTokens = @()
    foreach ($t in $tokens) {
        if ($t.ToLower() -eq "system") {
            $systemTokens += $t
        } else {
            $otherTokens += $t
        }
    }

    if ($systemTokens.Count -eq 0) {
        continue
    }

    $newTokens = $systemTokens + $otherTokens
    $newVal = $newTokens -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
